# Fix the typo in the document title: "Internview" -> "Interview".
# Word records the caret position of the last edit in the hidden
# "_GoBack" bookmark, so after this correction that bookmark is
# re-anchored at the point of the edit (inside the title run) instead
# of its old location further down in the document.

$d = $word.ActiveDocument

# 1) Correct the typo in the title.
$d.Content.Find.Execute(
    "Internview", $false, $false, $false, $false, $false,
    $true, 1, $false, "Interview", 2
)

# 2) Re-anchor the "_GoBack" bookmark at the edit point (right after
#    "Cracking the Coding Inter", before "view  - Notes"). Adding a
#    bookmark with this reserved name automatically removes it from
#    wherever it previously lived in the document.
$goBackRange = $d.Range(25, 25)
$d.Bookmarks.Add("_GoBack", $goBackRange)
